# Re-brand the "Product_KPI_Dashboard" template from the AI/ML sample
# content to generic Product Development content (15 Product templates
# with correct industry content).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Instructions & User Guide"
# ---------------------------------------------------------------------
$wsGuide = $wb.Worksheets.Item("Instructions & User Guide")

$wsGuide.Range("A1").Value = "Product Development KPI Dashboard - User Guide & Instructions"
$wsGuide.Range("B23").Value = "Availability and reliability of Product systems"

# ---------------------------------------------------------------------
# Sheet 2: "KPI Dashboard"
# ---------------------------------------------------------------------
$wsKpi = $wb.Worksheets.Item("KPI Dashboard")

$wsKpi.Range("A1").Value = "PRODUCT DEVELOPMENT - KPI DASHBOARD"
$wsKpi.Range("A2").Value = "Project: Product Development Implementation"

# Owner column ("Owner" = I) mentions of "ML Engineers" -> "Product Engineers"
$wsKpi.Range("I10").Value = "Product Engineers"
$wsKpi.Range("I16").Value = "Product Engineers"
$wsKpi.Range("I22").Value = "Product Engineers"

# Notes column (K) for every KPI row 8-22: swap the AI/ML blurb for the
# Product Development one.
for ($row = 8; $row -le 22; $row++) {
    $cell = "K" + $row
    $wsKpi.Range($cell).Value = "Critical KPI for Product Development success"
}
